# [IMP] New data for test environment
#
# Update the journal "code"/"name" shared strings on the account_journal
# sheet to point at the external.* identifiers instead of the old
# z0bug.jou_* ones, and shorten the sales-journal code from FATT to FAT.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Operazioni varie (misc) journal: id
$ws.Range("A7").Value = "external.MISC|VARIE"

# Fatture di vendita (sales) journal: id + code
$ws.Range("A8").Value = "external.FAT|FATT|INV"
$ws.Range("G8").Value = "FAT"

# Fatture di acquisto (purchase) journal: id
$ws.Range("A9").Value = "external.ACQ|FATTU|BILL"

# Column A grew a bit wider to fit the new, longer ids.
$ws.Columns.Item(1).ColumnWidth = 21.5

# Cursor/selection ends up parked on A7 after the edit.
$ws.Range("A7").Select() | Out-Null
